{"js": "// Replace the 15 lattice-multiplication problems in the (only) table with a\n// new set of problems, keeping the table's row/column layout and the\n// per-cell run formatting (sz=32) untouched.\n//\n// Each cell holds a single run whose text is 5 lines (separated by\n// manual line breaks):\n//   1) \"AB x CD\"            -- the two 2-digit factors\n//   2) \"  C    D\"           -- second factor's digits, spaced out\n//   3) \"  ----\"             -- a divider\n//   4) \"A|    |\"            -- first factor's tens digit\n//   5) \"B|    |\"            -- first factor's units digit\n\nconst newProblems = [\n  \"61 x 21\", \"74 x 38\", \"27 x 45\",\n  \"54 x 27\", \"37 x 76\", \"82 x 72\",\n  \"64 x 80\", \"89 x 95\", \"43 x 48\",\n  \"91 x 78\", \"45 x 90\", \"46 x 16\",\n  \"32 x 74\", \"47 x 71\", \"20 x 56\",\n];\n\nfunction buildCellText(problem) {\n  const m = problem.match(/^(\\d)(\\d) x (\\d)(\\d)$/);\n  const [, a, b, c, d] = m;\n  const LB = \"\\v\"; // Word manual line break (maps to <w:br/>)\n  return (\n    problem + LB +\n    \"  \" + c + \"    \" + d + LB +\n    \"  ----\" + LB +\n    a + \"|    |\" + LB +\n    b + \"|    |\"\n  );\n}\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nlet idx = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < 3; c++) {\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    para.insertText(buildCellText(newProblems[idx]), Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 15 lattice-multiplication problems in the (only) table with a\n# new set of problems, keeping the table's row/column layout and the\n# per-cell run formatting (sz=32) untouched.\n#\n# Each cell holds 5 lines (separated by manual line breaks, chr(11)):\n#   1) \"AB x CD\"            -- the two 2-digit factors\n#   2) \"  C    D\"           -- second factor's digits, spaced out\n#   3) \"  ----\"             -- a divider\n#   4) \"A|    |\"            -- first factor's tens digit\n#   5) \"B|    |\"            -- first factor's units digit\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$newProblems = @(\n  \"61 x 21\", \"74 x 38\", \"27 x 45\",\n  \"54 x 27\", \"37 x 76\", \"82 x 72\",\n  \"64 x 80\", \"89 x 95\", \"43 x 48\",\n  \"91 x 78\", \"45 x 90\", \"46 x 16\",\n  \"32 x 74\", \"47 x 71\", \"20 x 56\"\n)\n\n$LB = [char]11\n\n$rows = $tbl.Rows.Count\n$cols = $tbl.Columns.Count\n\n$idx = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    $problem = $newProblems[$idx]\n    $a = $problem.Substring(0,1)\n    $b = $problem.Substring(1,1)\n    $cc = $problem.Substring(5,1)\n    $dd = $problem.Substring(6,1)\n\n    $newText = $problem + $LB + \"  \" + $cc + \"    \" + $dd + $LB + \"  ----\" + $LB + $a + \"|    |\" + $LB + $b + \"|    |\"\n\n    $cell = $tbl.Cell($r, $c)\n    $cell.Range.Text = $newText\n\n    $idx++\n  }\n}\n"}
